$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in D/E are stored as text (matches original inlineStr cells)
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "59.928.00"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").Value = "2.556.27"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "502.98"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").Value = "152.25"
$ws.Range("E6").Value = "  -4.28%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  -6.52%  "
$ws.Range("D9").Value = "2.570.83"
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  +6.06%  "
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "3.013.11"
$ws.Range("E14").Value = "  +4.20%  "
$ws.Range("D15").Value = "60.029.95"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "21.45"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "2.574.28"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "4.78"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "345.23"
$ws.Range("E20").Value = "  +4.64%  "
$ws.Range("D21").Value = "10.21"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "6.01"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "0.995"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "59.85"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("D25").Value = "0.417"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "2.686.30"
$ws.Range("E27").Value = "  +4.61%  "
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "0.0₃0840"
$ws.Range("E29").Value = "  +3.97%  "
$ws.Range("D30").Value = "7.41"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "155.23"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").Value = "19.12"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").Value = "3.96"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  +20.34%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "3.75"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "0.839"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").Value = "1.45"
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").Value = "296.81"
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("D43").Value = "35.41"
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").Value = "0.0560"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0992"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.612"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "19.48"
$ws.Range("E48").Value = "  +7.51%  "
$ws.Range("D49").Value = "4.84"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").Value = "0.0232"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "10.27"
$ws.Range("E51").Value = "  -0.12%  "

# Restore default style (remove the temporary text-format style index) to match original formatting
$textRange.Style = "Normal"
